$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 19 (row 28) - CSK vs RCB
$ws.Range("E28").Value = 100
$ws.Range("H28").Value = 80
$ws.Range("K28").Value = 60
$ws.Range("N28").Value = 40
$ws.Range("Q28").Value = 0
$ws.Range("T28").Value = 20

# Contest 20 (row 29) - SRH vs DC
$ws.Range("E29").Value = 20
$ws.Range("H29").Value = 40
$ws.Range("K29").Value = 100
$ws.Range("N29").Value = 60
$ws.Range("Q29").Value = 0
$ws.Range("T29").Value = 80
